# Auto-applied cell updates for cryptos.xlsx (Price / Volume(1h) columns)
# Generated from the authoritative OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number need to be forced back to
# text (matching the original inlineStr cell type) without leaving a
# lasting style change on the cell: stamp Text format, write the
# value, then restore the Normal style so no 's' attribute sticks.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '60.552.16'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '2.898.49'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue 'D5' '527.76'
$ws.Range('E5').Value = '  -2.53%  '
Set-TextValue 'D6' '143.39'
$ws.Range('E6').Value = '  -5.81%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue 'D8' '0.555'
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').Value = '2.904.64'
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('E10').Value = '  -4.38%  '
Set-TextValue 'D11' '6.04'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '3.394.38'
$ws.Range('E13').Value = '  -2.40%  '
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').Value = '60.521.23'
$ws.Range('E15').Value = '  -1.73%  '
Set-TextValue 'D16' '22.76'
$ws.Range('E16').Value = '  -4.19%  '
$ws.Range('D17').Value = '2.903.04'
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('E18').Value = '  -3.92%  '
Set-TextValue 'D19' '5.04'
$ws.Range('E19').Value = '  -2.19%  '
Set-TextValue 'D20' '11.71'
$ws.Range('E20').Value = '  -2.24%  '
Set-TextValue 'D21' '361.19'
$ws.Range('E21').Value = '  -5.26%  '
Set-TextValue 'D22' '6.62'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('E23').Value = '  -0.03%  '
Set-TextValue 'D24' '5.69'
$ws.Range('E24').Value = '  +0.38%  '
Set-TextValue 'D25' '64.60'
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('E27').Value = '  -2.71%  '
Set-TextValue 'D28' '0.998'
$ws.Range('E28').Value = '  +0.12%  '
Set-TextValue 'D29' '7.87'
$ws.Range('E29').Value = '  -5.84%  '
$ws.Range('D30').Value = '0.0₃0852'
$ws.Range('E30').Value = '  -8.75%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -2.08%  '
Set-TextValue 'D33' '19.74'
$ws.Range('E33').Value = '  -3.22%  '
Set-TextValue 'D34' '152.44'
$ws.Range('E34').Value = '  -4.15%  '
$ws.Range('E35').Value = '  -6.23%  '
Set-TextValue 'D36' '5.59'
$ws.Range('E36').Value = '  -6.16%  '
$ws.Range('E37').Value = '  -5.79%  '
$ws.Range('E38').Value = '  -5.50%  '
Set-TextValue 'D39' '37.69'
$ws.Range('E39').Value = '  +1.48%  '
$ws.Range('E40').Value = '  -4.40%  '
Set-TextValue 'D41' '3.71'
$ws.Range('E41').Value = '  -5.66%  '
$ws.Range('D42').Value = '2.293.83'
$ws.Range('E42').Value = '  -4.71%  '
$ws.Range('E43').Value = '  -2.43%  '
Set-TextValue 'D44' '0.0582'
$ws.Range('E44').Value = '  -2.20%  '
Set-TextValue 'D45' '20.43'
$ws.Range('E45').Value = '  -7.90%  '
Set-TextValue 'D46' '0.997'
$ws.Range('E46').Value = '  +0.00%  '
Set-TextValue 'D47' '4.96'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('E48').Value = '  -3.19%  '
Set-TextValue 'D49' '10.31'
$ws.Range('E49').Value = '  -1.40%  '
Set-TextValue 'D50' '0.0923'
$ws.Range('E50').Value = '  -3.24%  '
Set-TextValue 'D51' '250.77'
$ws.Range('E51').Value = '  -6.33%  '
